$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at 895, pushing existing rows 895:929 down to 899:933
$ws.Rows("895:898").Insert()

# New row 895: Lechuga / Conconina(o) / Primera
$ws.Range("A895").Value = 7
$ws.Range("B895").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C895").Value = "Ñuble"
$ws.Range("D895").Value = 44939
$ws.Range("E895").Value = 16
$ws.Range("F895").Value = 100112033
$ws.Range("G895").Value = "Lechuga"
$ws.Range("H895").Value = "Conconina(o)"
$ws.Range("I895").Value = "Primera"
$ws.Range("J895").Value = 120
$ws.Range("K895").Value = 5000
$ws.Range("L895").Value = 5500
$ws.Range("M895").Value = 5250
$ws.Range("N895").Value = "$/caja 10 unidades"
$ws.Range("O895").Value = "Región del Maule"
$ws.Range("P895").Value = 525
$ws.Range("Q895").Value = 10
$ws.Range("R895").Value = "Hortaliza"

# New row 896: Lechuga / Conconina(o) / Primera
$ws.Range("A896").Value = 7
$ws.Range("B896").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C896").Value = "Ñuble"
$ws.Range("D896").Value = 44939
$ws.Range("E896").Value = 16
$ws.Range("F896").Value = 100112033
$ws.Range("G896").Value = "Lechuga"
$ws.Range("H896").Value = "Conconina(o)"
$ws.Range("I896").Value = "Primera"
$ws.Range("J896").Value = 120
$ws.Range("K896").Value = 6000
$ws.Range("L896").Value = 6500
$ws.Range("M896").Value = 6250
$ws.Range("N896").Value = "$/caja 15 unidades"
$ws.Range("O896").Value = "Región del Maule"
$ws.Range("P896").Value = 417
$ws.Range("Q896").Value = 15
$ws.Range("R896").Value = "Hortaliza"

# New row 897: Lechuga / Española / Primera
$ws.Range("A897").Value = 7
$ws.Range("B897").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C897").Value = "Ñuble"
$ws.Range("D897").Value = 44939
$ws.Range("E897").Value = 16
$ws.Range("F897").Value = 100112033
$ws.Range("G897").Value = "Lechuga"
$ws.Range("H897").Value = "Española"
$ws.Range("I897").Value = "Primera"
$ws.Range("J897").Value = 120
$ws.Range("K897").Value = 5000
$ws.Range("L897").Value = 5500
$ws.Range("M897").Value = 5250
$ws.Range("N897").Value = "$/caja 18 unidades"
$ws.Range("O897").Value = "Región del Maule"
$ws.Range("P897").Value = 292
$ws.Range("Q897").Value = 18
$ws.Range("R897").Value = "Hortaliza"

# New row 898: Lechuga / Marina / Primera
$ws.Range("A898").Value = 7
$ws.Range("B898").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C898").Value = "Ñuble"
$ws.Range("D898").Value = 44939
$ws.Range("E898").Value = 16
$ws.Range("F898").Value = 100112033
$ws.Range("G898").Value = "Lechuga"
$ws.Range("H898").Value = "Marina"
$ws.Range("I898").Value = "Primera"
$ws.Range("J898").Value = 120
$ws.Range("K898").Value = 5000
$ws.Range("L898").Value = 5500
$ws.Range("M898").Value = 5250
$ws.Range("N898").Value = "$/caja 18 unidades"
$ws.Range("O898").Value = "Región del Maule"
$ws.Range("P898").Value = 292
$ws.Range("Q898").Value = 18
$ws.Range("R898").Value = "Hortaliza"
